$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($CellRef, $Text) {
    $cell = $ws.Range($CellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = $origStyle
}

Set-TextValue "D2" "29.469.34"
Set-TextValue "E2" "  +0.15%  "
Set-TextValue "D3" "1.854.73"
Set-TextValue "E3" "  +0.32%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.29%  "
Set-TextValue "D5" "241.46"
Set-TextValue "E5" "  +0.40%  "
Set-TextValue "D6" "0.6330"
Set-TextValue "E6" "  +1.02%  "
Set-TextValue "D7" "1.003"
Set-TextValue "E7" "  +0.23%  "
Set-TextValue "D8" "0.07541"
Set-TextValue "E8" "  -2.00%  "
Set-TextValue "D9" "0.2923"
Set-TextValue "E9" "  +0.30%  "
Set-TextValue "D10" "24.52"
Set-TextValue "E10" "  -1.06%  "
Set-TextValue "D11" "0.07762"
Set-TextValue "E11" "  +0.20%  "
Set-TextValue "D12" "1.855.70"
Set-TextValue "E12" "  +0.71%  "
Set-TextValue "D13" "5.035"
Set-TextValue "E13" "  +0.22%  "
Set-TextValue "D14" "0.6838"
Set-TextValue "E14" "  +0.38%  "
Set-TextValue "D15" "0.00001043"
Set-TextValue "E15" "  -3.02%  "
Set-TextValue "D16" "83.37"
Set-TextValue "E16" "  -0.18%  "
Set-TextValue "D17" "2.120.79"
Set-TextValue "E17" "  +1.21%  "
Set-TextValue "D18" "6.143"
Set-TextValue "E18" "  -0.33%  "
Set-TextValue "D19" "29.502.66"
Set-TextValue "E19" "  +0.19%  "
Set-TextValue "D20" "229.78"
Set-TextValue "E20" "  +0.61%  "
Set-TextValue "D21" "12.39"
Set-TextValue "E21" "  +0.08%  "
Set-TextValue "D22" "1.003"
Set-TextValue "E22" "  +0.24%  "
Set-TextValue "D23" "7.496"
Set-TextValue "E23" "  +1.17%  "
Set-TextValue "D24" "1.003"
Set-TextValue "E24" "  +0.22%  "
Set-TextValue "D25" "159.47"
Set-TextValue "E25" "  +1.55%  "
Set-TextValue "D26" "0.1390"
Set-TextValue "E26" "  +1.18%  "
Set-TextValue "D27" "8.457"
Set-TextValue "E27" "  +0.74%  "
Set-TextValue "D28" "17.66"
Set-TextValue "E28" "  -0.18%  "
Set-TextValue "D29" "1.430"
Set-TextValue "E29" "  +6.62%  "
Set-TextValue "D30" "1.479"
Set-TextValue "E30" "  +1.08%  "
Set-TextValue "D31" "0.05714"
Set-TextValue "E31" "  +1.26%  "
Set-TextValue "D32" "4.147"
Set-TextValue "E32" "  +0.77%  "
Set-TextValue "D33" "4.054"
Set-TextValue "E33" "  +0.53%  "
Set-TextValue "D34" "1.158"
Set-TextValue "E34" "  -0.21%  "
Set-TextValue "D35" "1.824"
Set-TextValue "E35" "  -0.92%  "
Set-TextValue "D36" "0.6960"
Set-TextValue "E36" "  -1.67%  "
Set-TextValue "D37" "2.595"
Set-TextValue "E37" "  -0.07%  "
Set-TextValue "B38" "MXToken"
Set-TextValue "C38" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D38" "2.828"
Set-TextValue "E38" "  +2.22%  "
Set-TextValue "B39" "Maker"
Set-TextValue "C39" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D39" "1.254.02"
Set-TextValue "E39" "  +2.11%  "
Set-TextValue "D40" "0.01831"
Set-TextValue "E40" "  +2.36%  "
Set-TextValue "D41" "6.494"
Set-TextValue "E41" "  +0.46%  "
Set-TextValue "D42" "0.9081"
Set-TextValue "E42" "  +0.30%  "
Set-TextValue "D43" "1.003"
Set-TextValue "E43" "  +0.16%  "
Set-TextValue "D44" "2.021.92"
Set-TextValue "E44" "  +0.95%  "
Set-TextValue "D45" "101.61"
Set-TextValue "E45" "  -0.02%  "
Set-TextValue "D46" "66.14"
Set-TextValue "E46" "  +0.46%  "
Set-TextValue "D47" "7.120"
Set-TextValue "E47" "  -0.45%  "
Set-TextValue "D48" "0.1168"
Set-TextValue "E48" "  +0.96%  "
Set-TextValue "B49" "BabyDogeCoin"
Set-TextValue "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.00000000117"
Set-TextValue "E49" "  -2.77%  "
Set-TextValue "B50" "EnergySwap"
Set-TextValue "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "9.068"
Set-TextValue "E50" "  +0.64%  "
Set-TextValue "B51" "TheSandbox"
Set-TextValue "C51" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D51" "0.3968"
Set-TextValue "E51" "  -1.01%  "

Write-Output "All cells updated."
